# Update Slit3-Robo2 NATMI LR-pair data with recomputed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.412972666666667
$ws.Range("H2").Value = 7.238918
$ws.Range("I2").Value = 0.04202116219916398
$ws.Range("J2").Value = 0.04202116219916398
$ws.Range("M2").Value = 0.2313123333333333
$ws.Range("N2").Value = 0.693937
$ws.Range("O2").Value = 0.7569517164947553
$ws.Range("P2").Value = 0.7569517164947555
$ws.Range("Q2").Value = 0.5581503377962223
$ws.Range("R2").Value = 5.023353040166
$ws.Range("S2").Value = 0.0318079908557617
$ws.Range("T2").Value = 0.03180799085576171

# Row 3
$ws.Range("G3").Value = 2.412972666666667
$ws.Range("H3").Value = 7.238918
$ws.Range("I3").Value = 0.04202116219916398
$ws.Range("J3").Value = 0.04202116219916398
$ws.Range("O3").Value = 0.2385552472206224
$ws.Range("P3").Value = 0.2385552472206224
$ws.Range("Q3").Value = 0.1759024901031111
$ws.Range("R3").Value = 1.583122410928
$ws.Range("S3").Value = 0.01002436873691943
$ws.Range("T3").Value = 0.01002436873691943

# Row 4
$ws.Range("G4").Value = 2.412972666666667
$ws.Range("H4").Value = 7.238918
$ws.Range("I4").Value = 0.04202116219916398
$ws.Range("J4").Value = 0.04202116219916398
$ws.Range("M4").Value = 0.001373
$ws.Range("N4").Value = 0.004119
$ws.Range("O4").Value = 0.004493036284622232
$ws.Range("P4").Value = 0.004493036284622232
$ws.Range("Q4").Value = 0.003313011471333334
$ws.Range("R4").Value = 0.029817103242
$ws.Range("S4").Value = 0.0001888026064828399
$ws.Range("T4").Value = 0.0001888026064828399

# Row 5
$ws.Range("I5").Value = 0.8969165968468352
$ws.Range("J5").Value = 0.8969165968468353
$ws.Range("M5").Value = 0.2313123333333333
$ws.Range("N5").Value = 0.693937
$ws.Range("O5").Value = 0.7569517164947553
$ws.Range("P5").Value = 0.7569517164947555
$ws.Range("Q5").Value = 11.91338543023589
$ws.Range("R5").Value = 107.220468872123
$ws.Range("S5").Value = 0.6789225575358464
$ws.Range("T5").Value = 0.6789225575358466

# Row 6
$ws.Range("I6").Value = 0.8969165968468352
$ws.Range("J6").Value = 0.8969165968468353
$ws.Range("O6").Value = 0.2385552472206224
$ws.Range("P6").Value = 0.2385552472206224
$ws.Range("S6").Value = 0.213964160497076
$ws.Range("T6").Value = 0.2139641604970761

# Row 7
$ws.Range("I7").Value = 0.8969165968468352
$ws.Range("J7").Value = 0.8969165968468353
$ws.Range("M7").Value = 0.001373
$ws.Range("N7").Value = 0.004119
$ws.Range("O7").Value = 0.004493036284622232
$ws.Range("P7").Value = 0.004493036284622232
$ws.Range("Q7").Value = 0.07071425012233333
$ws.Range("R7").Value = 0.6364282511010001
$ws.Range("S7").Value = 0.004029878813912721
$ws.Range("T7").Value = 0.004029878813912721

# Row 8
$ws.Range("G8").Value = 0.003190333333333333
$ws.Range("H8").Value = 0.009571
$ws.Range("I8").Value = 0.00005555865440224608
$ws.Range("J8").Value = 0.00005555865440224608
$ws.Range("M8").Value = 0.2313123333333333
$ws.Range("N8").Value = 0.693937
$ws.Range("O8").Value = 0.7569517164947553
$ws.Range("P8").Value = 0.7569517164947555
$ws.Range("Q8").Value = 0.0007379634474444445
$ws.Range("R8").Value = 0.006641671027
$ws.Range("S8").Value = 0.00004205521881591907
$ws.Range("T8").Value = 0.00004205521881591908

# Row 9
$ws.Range("G9").Value = 0.003190333333333333
$ws.Range("H9").Value = 0.009571
$ws.Range("I9").Value = 0.00005555865440224608
$ws.Range("J9").Value = 0.00005555865440224608
$ws.Range("O9").Value = 0.2385552472206224
$ws.Range("P9").Value = 0.2385552472206224
$ws.Range("Q9").Value = 0.0002325710462222222
$ws.Range("R9").Value = 0.002093139416
$ws.Range("S9").Value = 0.00001325380853617293
$ws.Range("T9").Value = 0.00001325380853617293

# Row 10
$ws.Range("G10").Value = 0.003190333333333333
$ws.Range("H10").Value = 0.009571
$ws.Range("I10").Value = 0.00005555865440224608
$ws.Range("J10").Value = 0.00005555865440224608
$ws.Range("M10").Value = 0.001373
$ws.Range("N10").Value = 0.004119
$ws.Range("O10").Value = 0.004493036284622232
$ws.Range("P10").Value = 0.004493036284622232
$ws.Range("Q10").Value = 0.000004380327666666667
$ws.Range("R10").Value = 0.000039422949
$ws.Range("S10").Value = 0.0000002496270501540783
$ws.Range("T10").Value = 0.0000002496270501540783

# Row 11
$ws.Range("G11").Value = 3.460560666666667
$ws.Range("H11").Value = 10.381682
$ws.Range("I11").Value = 0.06026457865970317
$ws.Range("J11").Value = 0.06026457865970317
$ws.Range("M11").Value = 0.2313123333333333
$ws.Range("N11").Value = 0.693937
$ws.Range("O11").Value = 0.7569517164947553
$ws.Range("P11").Value = 0.7569517164947555
$ws.Range("Q11").Value = 0.8004703624482224
$ws.Range("R11").Value = 7.204233262034001
$ws.Range("S11").Value = 0.04561737626029552
$ws.Range("T11").Value = 0.04561737626029552

# Row 12
$ws.Range("G12").Value = 3.460560666666667
$ws.Range("H12").Value = 10.381682
$ws.Range("I12").Value = 0.06026457865970317
$ws.Range("J12").Value = 0.06026457865970317
$ws.Range("O12").Value = 0.2385552472206224
$ws.Range("P12").Value = 0.2385552472206224
$ws.Range("Q12").Value = 0.2522702585191112
$ws.Range("R12").Value = 2.270432326672
$ws.Range("S12").Value = 0.01437643146081213
$ws.Range("T12").Value = 0.01437643146081213

# Row 13
$ws.Range("G13").Value = 3.460560666666667
$ws.Range("H13").Value = 10.381682
$ws.Range("I13").Value = 0.06026457865970317
$ws.Range("J13").Value = 0.06026457865970317
$ws.Range("M13").Value = 0.001373
$ws.Range("N13").Value = 0.004119
$ws.Range("O13").Value = 0.004493036284622232
$ws.Range("P13").Value = 0.004493036284622232
$ws.Range("Q13").Value = 0.004751349795333334
$ws.Range("R13").Value = 0.04276214815800001
$ws.Range("S13").Value = 0.0002707709385955169
$ws.Range("T13").Value = 0.0002707709385955169

# Row 14
$ws.Range("G14").Value = 0.03242133333333333
$ws.Range("H14").Value = 0.09726399999999999
$ws.Range("I14").Value = 0.0005646073515599271
$ws.Range("J14").Value = 0.0005646073515599271
$ws.Range("M14").Value = 0.2313123333333333
$ws.Range("N14").Value = 0.693937
$ws.Range("O14").Value = 0.7569517164947553
$ws.Range("P14").Value = 0.7569517164947555
$ws.Range("Q14").Value = 0.00749945426311111
$ws.Range("R14").Value = 0.067495088368
$ws.Range("S14").Value = 0.0004273805039088446
$ws.Range("T14").Value = 0.0004273805039088447

# Row 15
$ws.Range("G15").Value = 0.03242133333333333
$ws.Range("H15").Value = 0.09726399999999999
$ws.Range("I15").Value = 0.0005646073515599271
$ws.Range("J15").Value = 0.0005646073515599271
$ws.Range("O15").Value = 0.2385552472206224
$ws.Range("P15").Value = 0.2385552472206224
$ws.Range("Q15").Value = 0.002363471971555555
$ws.Range("R15").Value = 0.021271247744
$ws.Range("S15").Value = 0.0001346900463339592
$ws.Range("T15").Value = 0.0001346900463339593

# Row 16
$ws.Range("G16").Value = 0.03242133333333333
$ws.Range("H16").Value = 0.09726399999999999
$ws.Range("I16").Value = 0.0005646073515599271
$ws.Range("J16").Value = 0.0005646073515599271
$ws.Range("M16").Value = 0.001373
$ws.Range("N16").Value = 0.004119
$ws.Range("O16").Value = 0.004493036284622232
$ws.Range("P16").Value = 0.004493036284622232
$ws.Range("Q16").Value = 0.00004451449066666667
$ws.Range("R16").Value = 0.000400630416
$ws.Range("S16").Value = 0.000002536801317123213
$ws.Range("T16").Value = 0.000002536801317123213

# Row 17
$ws.Range("G17").Value = 0.01019233333333333
$ws.Range("H17").Value = 0.030577
$ws.Range("I17").Value = 0.0001774962883353337
$ws.Range("J17").Value = 0.0001774962883353337
$ws.Range("M17").Value = 0.2313123333333333
$ws.Range("N17").Value = 0.693937
$ws.Range("O17").Value = 0.7569517164947553
$ws.Range("P17").Value = 0.7569517164947555
$ws.Range("Q17").Value = 0.002357612405444444
$ws.Range("R17").Value = 0.021218511649
$ws.Range("S17").Value = 0.0001343561201268788
$ws.Range("T17").Value = 0.0001343561201268788

# Row 18
$ws.Range("G18").Value = 0.01019233333333333
$ws.Range("H18").Value = 0.030577
$ws.Range("I18").Value = 0.0001774962883353337
$ws.Range("J18").Value = 0.0001774962883353337
$ws.Range("O18").Value = 0.2385552472206224
$ws.Range("P18").Value = 0.2385552472206224
$ws.Range("Q18").Value = 0.0007430075102222222
$ws.Range("R18").Value = 0.006687067592
$ws.Range("S18").Value = 0.00004234267094457839
$ws.Range("T18").Value = 0.00004234267094457839

# Row 19
$ws.Range("G19").Value = 0.01019233333333333
$ws.Range("H19").Value = 0.030577
$ws.Range("I19").Value = 0.0001774962883353337
$ws.Range("J19").Value = 0.0001774962883353337
$ws.Range("M19").Value = 0.001373
$ws.Range("N19").Value = 0.004119
$ws.Range("O19").Value = 0.004493036284622232
$ws.Range("P19").Value = 0.004493036284622232
$ws.Range("Q19").Value = 0.00001399407366666667
$ws.Range("R19").Value = 0.000125946663
$ws.Range("S19").Value = 0.0000007974972638764239
$ws.Range("T19").Value = 0.0000007974972638764239
